$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 3 new rows above the "MAX(...)+1" test-counter row (old row 70) ---
# This pushes the old row 70 down to row 73, and every row after it down by 3
# (old 72 -> 75, old 73..94 -> 76..97, old 96 -> 99, old 97 -> 100), exactly
# matching the target layout.
$ws.Rows("70:72").Insert()

# --- Row 69 (existing blank row, now holding the new "Gist URL:" label) ---
# Re-use the formatting already used for the similar label in D66
# ("Minimised version:") so no stray new style is introduced.
$ws.Range("D66").Copy()
$ws.Range("D69").PasteSpecial(-4122)
$ws.Range("D69").Value = "Gist URL:"

# --- Row 70 (first of the freshly inserted rows) holds the hyperlink ---
# Give E/F on both new rows the same "filler" formatting used elsewhere in
# this section (e.g. E20:F20) before the hyperlink creates its own style.
$ws.Range("E20:F20").Copy()
$ws.Range("E69:F69").PasteSpecial(-4122)
$ws.Range("E20:F20").Copy()
$ws.Range("E70:F70").PasteSpecial(-4122)

$gistUrl = "https://gist.github.com/r-silk/8ca0742d549ec4153d403b0a64847974"
$ws.Hyperlinks.Add($ws.Range("D70"), $gistUrl, "", "", $gistUrl)

# --- Fix up the "running error id" counter formula, now on row 73 ---
# It should keep referencing the blank row directly above the gap, which is
# now row 71 (was row 69 before the insert).
$ws.Range("B73").Formula = "=MAX(B`$1:B71) + 1"

$wb.Application.CutCopyMode = $false
